# 3 Mayis 2020 verileri eklendi
# Adds the 2020-05-03 COVID-19 Turkey data row to the "data" sheet and
# grows the worksheet Table (Table3) by one row so it keeps covering the
# full data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's data lives inside an Excel Table - use ListRows.Add() so the
# table range/autoFilter/dimension all expand together, the same as typing
# a new row directly below an existing table in the real app.
$lo = $ws.ListObjects.Item("Table3")
$lo.ListRows.Add()

# New row is row 53 (table header is row 1, so row 52 was the last data row).
$ws.Range("A53").Value = 43954
$ws.Range("B53").Value = 24001
$ws.Range("C53").Value = 1670
$ws.Range("D53").Value = 61
$ws.Range("E53").Value = 4892

# Leave the selection on the last cell entered, matching manual data entry.
$ws.Range("E53").Select()
